$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CasesTab query (B2): append an ORDER BY / LIMIT clause ---
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + "`norder By ss.study_subject_id ASC LIMIT 100 "

# --- SamplesTab query (B3): append an ORDER BY / LIMIT clause ---
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# --- FilesTab query (B4): replace the lower-case "order by" clause ---
$b4 = $ws.Range("B4").Value2
$b4 = $b4.Replace("    order by f.file_name", "    order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value2 = $b4

# --- Row heights grew because the wrapped text now spans one more line ---
$ws.Rows(2).RowHeight = 360
$ws.Rows(3).RowHeight = 374.4

# --- Window/selection moved down one row (scrolled so row 3 is at top, B4 selected) ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B4").Select()
